# Update simulation results: refresh the stats for generations 0-9 with the
# latest run's numbers, and drop the now-stale generations 10-14 (rows 12-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-generation data (Prey, Pred, Time(ms), prey starved, predator starved)
# for rows 2-11 (Generation 0-9). Column A (Generation) and B (Winner) are
# unchanged.
$data = @{
    2  = @(37, 20, 2199.702739715576, 10, 0)
    3  = @(38, 3,  1629.07862663269,  8,  17)
    4  = @(41, 5,  1659.323692321777, 4,  15)
    5  = @(43, 3,  1700.087308883667, 4,  17)
    6  = @(43, 3,  1690.253496170044, 4,  17)
    7  = @(44, 3,  1680.027484893799, 3,  17)
    8  = @(43, 3,  1672.077894210815, 4,  17)
    9  = @(43, 4,  1645.671606063843, 4,  16)
    10 = @(44, 4,  1636.59143447876,  3,  16)
    11 = @(43, 7,  1726.107835769653, 3,  13)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("C$r").Value = $vals[0]
    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
    $ws.Range("F$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}

# The simulation now only ran 10 generations (rows 2-11); remove the old
# trailing rows for generations 10-14 that no longer exist.
$ws.Rows("12:16").Delete()
